# Insert a new data row at row 111 (pushing the existing rows 111-231 down
# to 112-232) and populate the new row with the latest weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(111).Insert()

$ws.Range("A111").Value = 9
$ws.Range("B111").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C111").Value = "Metropolitana"
$ws.Range("D111").Value = 44629
$ws.Range("E111").Value = 13
$ws.Range("F111").Value = 300000001
$ws.Range("G111").Value = "Rabanito"
$ws.Range("H111").Value = "Sin especificar"
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 5200
$ws.Range("K111").Value = 3000
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 3000
$ws.Range("N111").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O111").Value = "Provincia de Chacabuco"
$ws.Range("P111").Value = 30
$ws.Range("Q111").Value = 100
$ws.Range("R111").Value = "Hortaliza"
